$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held a stray "Docentes responsáveis" value (no label in
# column A) is removed entirely; everything below it shifts up by one row,
# taking the sheet from 22 data rows down to 21 (dimension A1:C21).
$ws.Rows(13).Delete()

# Row 10 (Objetivos:) now carries the teacher's name instead of the old
# long objectives description.
$ws.Range("B10:C10").Value = "5840671 - Francisco José Moreira Chaves"

# Row 13 (Programa resumido:, shifted up from old row 14) now just says
# "Semestral" instead of the old long summary text.
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 (Programa:, shifted up from old row 16) now carries the activation
# date text instead of the old long program text. Assign it via a formula
# and flatten to a literal value so Excel's automatic date-pattern parsing
# doesn't turn the text into a date serial number (we need literal text,
# matching the rest of the sheet), while keeping the original cell style.
$ws.Range("B15:C15").Formula = '="01/01/2018"'
$ws.Range("B15:C15").Copy()
$ws.Range("B15:C15").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Row 18 (Método:, shifted up from old row 19) now carries the teacher's
# name value (same stray value as row 10).
$ws.Range("B18:C18").Value = "5840671 - Francisco José Moreira Chaves"

# Row 19 (Critério:, shifted up from old row 20) now carries the old
# "Método" evaluation text.
$ws.Range("B19:C19").Value = "Duas Provas  P1  1º bimestre e P2  2º bimestre"

# Row 20 (Norma de recuperação:, shifted up from old row 21) now carries the
# old "Critério" formula text.
$ws.Range("B20:C20").Value = "MF = (P1+ P2)/2"

# Row 21 (Bibliografia:, shifted up from old row 22) now carries the old
# "Norma de recuperação" formula text; the long bibliography text is gone.
$ws.Range("B21:C21").Value = "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"
